$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.346.02"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "1.862.17"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.63"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.32%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4766"
$ws.Range("D7").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2757"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06449"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.28%  "
$ws.Range("D10").Value = "1.863.66"
$ws.Range("E10").Value = "  -0.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07436"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("E12").Value = "  -2.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.990"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "85.81"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6325"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.69%  "
$ws.Range("D16").Value = "30.291.64"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.61"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.79"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007372"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.20%  "
$ws.Range("D21").Value = "2.096.83"
$ws.Range("E21").Value = "  -3.69%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  -4.08%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.010"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.13%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.284"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.53"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.89"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.26%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.863"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.63%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.382"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.42%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09993"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +5.02%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.217"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.922"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.72%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04919"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.17%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.148"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.44%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7232"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.41%  "
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9995"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.697"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01930"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.41%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.632"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9045"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.60%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.985"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.37%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.65"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4107"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.92%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.540"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.56%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.054"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.65%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "61.28"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.19%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1208"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.92%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.815"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.402"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.02%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.10"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.13%  "
